$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new person type row: ID 4 -> "retired non working adult"
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "retired non working adult"

# Match the selection noted in the saved file
$ws.Range("B6").Select()
